# "lógica para não gerar a mesma rota 2 vezes"
# Rebuild the delivery-route table (A2:E16) with the deduplicated/re-ordered
# routes. Rows 2-9 already exist and get their values corrected; rows 10-16
# are brand-new rows that are appended (formatted like the existing data
# rows) and the sheet's dimension grows from A1:E10 to A1:E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the data block with new rows 11-16, copying the look & feel
#     (borders/alignment) of the existing data rows --------------------
$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A10:E16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- column A holds order numbers stored as text (e.g. "11066"), so force
#     a text number format before writing them, otherwise Excel would
#     silently convert the literal into a numeric value -----------------
$ws.Range("A2:A16").NumberFormat = "@"

$rows = @(
    @(2,  "11066", "13/11/2024", "manhã", "Desconhecido", "Barra do Aririú"),
    @(3,  "11066", "13/11/2024", "manhã", "Desconhecido", "Barra do Aririú"),
    @(4,  "11065", "13/11/2024", "manhã", "Desconhecido", "Barra do Aririú"),
    @(5,  "11065", "13/11/2024", "manhã", "Desconhecido", "Barra do Aririú"),
    @(6,  "11104", "13/11/2024", "manhã", "Desconhecido", "Caminho Novo"),
    @(7,  "11104", "13/11/2024", "tarde", "Desconhecido", "Caminho Novo"),
    @(8,  "11024", "13/11/2024", "tarde", "Desconhecido", "Agronômica"),
    @(9,  "11024", "13/11/2024", "tarde", "Desconhecido", "Agronômica"),
    @(10, "11121", "13/11/2024", "tarde", "Desconhecido", "Canto"),
    @(11, "11103", "14/11/2024", "manhã", "Desconhecido", "Caminho Novo"),
    @(12, "11103", "14/11/2024", "manhã", "Desconhecido", "Caminho Novo"),
    @(13, "11058", "14/11/2024", "manhã", "Desconhecido", "Aririú"),
    @(14, "11058", "14/11/2024", "tarde", "Desconhecido", "Aririú"),
    @(15, "11060", "14/11/2024", "tarde", "Desconhecido", "Aririú"),
    @(16, "11060", "14/11/2024", "tarde", "Desconhecido", "Aririú")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
